$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column O (year 2021) that mirrors column N (year 2020),
# copying values and formatting, then apply the handful of cells whose
# values differ from 2020.
$ws.Range("N4:N14").Copy($ws.Range("O4:O14"))

# Header: year label
$ws.Range("O4").Value = 2021

# Row-by-row 2021 values (most copied unchanged from 2020, a few differ)
$ws.Range("O5").Value = 2
$ws.Range("O6").Value = "-"
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = "-"
$ws.Range("O9").Value = "-"
$ws.Range("O10").Value = "-"
$ws.Range("O11").Value = "-"
$ws.Range("O12").Value = 1
$ws.Range("O13").Value = "-"
$ws.Range("O14").Value = "-"

# Match the selection saved in the source file
$ws.Range("P1").Select()
